$d = $word.ActiveDocument

# --- Fix grammar in the "payForStay" bullet: "vehicles is" -> "vehicle is"
# and remove the duplicated "it" before "you should increase".
#
# Original: "Check if the vehicles is parked and if it is it you should increase"
# Target:   "Check if the vehicle is parked and if it is you should increase"
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Check if the vehicles is parked and if it is it you should increase",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Check if the vehicle is parked and if it is you should increase",
    2)
if (-not $found1) {
    Write-Host "WARNING: first target sentence not found"
}

# --- Fix grammar in the following sentence of the same bullet:
# "However, if the vehicles is not parked" -> "However, if the vehicle is not parked"
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "However, if the vehicles is not parked, throw a new error ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "However, if the vehicle is not parked, throw a new error ",
    2)
if (-not $found2) {
    Write-Host "WARNING: second target sentence not found"
}

Write-Host "Done. found1=$found1 found2=$found2"
